$d = $word.ActiveDocument

function Get-ParagraphByText($needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1. "EC7" + " - ENGENHARIA DE SOFTWARE 2" -> single run (no visible text
#    change, just a run merge).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("EC7 - ENGENHARIA DE SOFTWARE 2", $true, $false, $false, $false, $false, $true, 1, $false, "EC7 - ENGENHARIA DE SOFTWARE 2", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Question text rewritten; keep the trailing "?" as its own run.
# ---------------------------------------------------------------------------
$oldQuestion = "Quando a semana de 40 horas não é suficiente para entregar o Sprint, o que fazer"
$newQuestion = "Porque o XP define semanas de no máximo 40 horas"
$d.Content.Find.Execute($oldQuestion, $true, $false, $false, $false, $false, $true, 1, $false, $newQuestion, 2) | Out-Null

$qPara = Get-ParagraphByText($newQuestion)
$qRange = $qPara.Range
$qXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
  '<w:p><w:r><w:rPr><w:b/></w:rPr><w:t>' + $newQuestion + '</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>?</w:t></w:r></w:p>' +
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$qTextOnly = $d.Range($qRange.Start, $qRange.End - 1)
$qTextOnly.InsertXML($qXml)

# ---------------------------------------------------------------------------
# 3. Rewrite the "A metodologia do XP pre" paragraph entirely and append a
#    new empty paragraph with matching formatting right after it.
# ---------------------------------------------------------------------------
$incompletePara = Get-ParagraphByText("A metodologia do XP pre")
$incompleteRange = $incompletePara.Range

$rPrXml = '<w:rPr><w:i/><w:color w:val="1F3864" w:themeColor="accent5" w:themeShade="80"/></w:rPr>'
$bodyXml = '<w:p><w:pPr><w:tabs><w:tab w:val="right" w:pos="7938"/></w:tabs><w:ind w:firstLine="709"/>' + $rPrXml + '</w:pPr>' +
  '<w:r>' + $rPrXml + '<w:tab/></w:r>' +
  '<w:r>' + $rPrXml + '<w:t>A metodologia do XP propõe um limite de 40 horas semanais, se durante a Sprint for necessário fazer horas extras para entregar os artefatos, significa que a Sprint não foi bem planejada. E a prática de fazer horas extras pode tornar-se um vício para as próximas Sprints. Por esse motivo, as horas extras não são bem vistas pelo XP</w:t></w:r>' +
  '<w:r>' + $rPrXml + '<w:t xml:space="preserve"> e devem ser evitadas sempre que possível.</w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
  '</w:p>' +
  '<w:p><w:pPr><w:tabs><w:tab w:val="right" w:pos="7938"/></w:tabs><w:ind w:firstLine="709"/>' + $rPrXml + '</w:pPr></w:p>'
$fullXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$incompleteRange.InsertXML($fullXml)

# ---------------------------------------------------------------------------
# 4. Bump several "firstLine" paragraph indents from 567 twips (28.35pt) to
#    709 twips (35.45pt) on the four answer paragraphs that changed.
# ---------------------------------------------------------------------------
$indentTargets = @("técnica de desenvolvimento", "Comunicação, feedback", "Feedback rápido", "Jogos de planejamento")
foreach ($needle in $indentTargets) {
    $p = Get-ParagraphByText($needle)
    $p.Format.FirstLineIndent = 35.45
}

# ---------------------------------------------------------------------------
# 5. Merge runs that were split mid-word/mid-phrase back into single runs
#    (no visible text change).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Semana de 40 h", $true, $false, $false, $false, $false, $true, 1, $false, "Semana de 40 h", 2) | Out-Null
$d.Content.Find.Execute("Propriedade coletiva", $true, $false, $false, $false, $false, $true, 1, $false, "Propriedade coletiva", 2) | Out-Null
$d.Content.Find.Execute("Programação em pares", $true, $false, $false, $false, $false, $true, 1, $false, "Programação em pares", 2) | Out-Null
$d.Content.Find.Execute("Padrões de codificação", $true, $false, $false, $false, $false, $true, 1, $false, "Padrões de codificação", 2) | Out-Null
